$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename sheet
$ws.Name = "Expanded-High"

# New data: years 2031-2048 with updated cumulative capacity values
$years = @(2031, 2032, 2033, 2034, 2035, 2036, 2037, 2038, 2039, 2040, 2041, 2042, 2043, 2044, 2045, 2046, 2047, 2048)
$values = @(
    1014.686248331108,
    1723.928200563715,
    3458.369223514548,
    5347.234032355945,
    7467.944060857642,
    10249.45008680604,
    13085.75617424004,
    16109.44111882,
    19527.0433267734,
    23355.60088930804,
    27566.85963634356,
    32109.79270138988,
    36648.19071634963,
    41101.57988587655,
    45470.40220569386,
    49453.25267151307,
    52755.98331132648,
    54645.70030435042
)

for ($i = 0; $i -lt $years.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $years[$i]
    $ws.Cells.Item($row, 2).Value = $values[$i]
}

# Remove leftover rows 20-27 (previously rows for years 2049-2057)
$ws.Range("A20:B27").Clear()
